$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "28.568.47"
$ws.Range("E2").Value2 = "  +0.50%  "
$ws.Range("D3").Value2 = "1.577.79"
$ws.Range("E3").Value2 = "  -0.66%  "
$ws.Range("E4").Value2 = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "213.23"
$ws.Range("E5").Value2 = "  -0.12%  "
$ws.Range("E6").Value2 = "  -0.51%  "
$ws.Range("E8").Value2 = "  +1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "24.00"
$ws.Range("E9").Value2 = "  -1.63%  "
$ws.Range("E10").Value2 = "  -1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.0591"
$ws.Range("E11").Value2 = "  -1.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "0.0894"
$ws.Range("E12").Value2 = "  +0.77%  "
$ws.Range("D13").Value2 = "1.803.39"
$ws.Range("E13").Value2 = "  -0.64%  "
$ws.Range("D14").Value2 = "1.577.99"
$ws.Range("E14").Value2 = "  -1.01%  "
$ws.Range("D15").Value2 = "28.577.16"
$ws.Range("E16").Value2 = "  -1.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.518"
$ws.Range("E17").Value2 = "  -2.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "62.16"
$ws.Range("E18").Value2 = "  -1.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "231.02"
$ws.Range("E19").Value2 = "  +0.53%  "
$ws.Range("E20").Value2 = "  -1.15%  "
$ws.Range("D21").Value2 = "0.0₃0691"
$ws.Range("E21").Value2 = "  -2.18%  "
$ws.Range("E22").Value2 = "  +0.09%  "
$ws.Range("E23").Value2 = "  -4.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "9.14"
$ws.Range("E24").Value2 = "  -2.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "2.07"
$ws.Range("E25").Value2 = "  +5.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "151.11"
$ws.Range("E26").Value2 = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "15.01"
$ws.Range("E27").Value2 = "  -1.29%  "
$ws.Range("E28").Value2 = "  -2.06%  "
$ws.Range("E29").Value2 = "  -2.55%  "
$ws.Range("E30").Value2 = "  +0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.0483"
$ws.Range("E31").Value2 = "  +2.38%  "
$ws.Range("E32").Value2 = "  -2.33%  "
$ws.Range("E33").Value2 = "  -1.30%  "
$ws.Range("E34").Value2 = "  -2.28%  "
$ws.Range("D35").Value2 = "1.398.65"
$ws.Range("E35").Value2 = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "1.06"
$ws.Range("E36").Value2 = "  +4.82%  "
$ws.Range("E38").Value2 = "  +0.61%  "
$ws.Range("E39").Value2 = "  +2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.0165"
$ws.Range("E40").Value2 = "  -0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.521"
$ws.Range("E41").Value2 = "  -3.78%  "
$ws.Range("E42").Value2 = "  +0.10%  "
$ws.Range("E43").Value2 = "  -2.16%  "
$ws.Range("E44").Value2 = "  +0.40%  "
$ws.Range("E45").Value2 = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "5.46"
$ws.Range("E46").Value2 = "  -1.84%  "
$ws.Range("E47").Value2 = "  -1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "63.06"
$ws.Range("E48").Value2 = "  -1.69%  "
$ws.Range("D49").Value2 = "1.715.30"
$ws.Range("E49").Value2 = "  -0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "86.56"
$ws.Range("E50").Value2 = "  -0.70%  "
$ws.Range("E51").Value2 = "  -1.49%  "
